$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 'AU-5 b,AU-5 a'
$ws.Range("A13").Value = 'AU-12 a,AU-3,CM-5 (1),MA-4 (1) (a),CM-6 b,AU-6 (4),AU-7 (1),AU-14 (1),AU-7 a,AU-3 (1)'
$ws.Range("A17").Value = 'CM-7 (2),CM-6 b'
$ws.Range("A22").Value = 'CM-7 (2),CM-6 b'
$ws.Range("A23").Value = 'CM-7 (2),CM-6 b'
$ws.Range("A37").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-14 (1),AU-12 c,AU-3 (1)'
$ws.Range("A38").Value = 'AU-14 (1),AU-4'
$ws.Range("A39").Value = 'AU-3,AU-4 (1)'
$ws.Range("A44").Value = 'AU-6 (4),CM-6 b,AU-4 (1)'
$ws.Range("A48").Value = 'IA-2 (11),IA-2 (12)'
$ws.Range("A49").Value = 'IA-2 (11),IA-2 (1),IA-2 (12)'
$ws.Range("A50").Value = 'SI-6 d,SI-6 b,CM-3 (5)'
$ws.Range("A51").Value = 'SI-6 d,CM-3 (5)'
$ws.Range("A52").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A53").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A54").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A55").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A56").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A57").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A58").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A59").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A60").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A61").Value = 'AU-12 a,CM-5 (1),CM-6 b,AU-8 b,AU-7 a,AU-7 b,AU-12 c,AU-12 (3)'
$ws.Range("A62").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A63").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A64").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A65").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A66").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A67").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A68").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A69").Value = 'SI-6 a,CM-3 (5)'
$ws.Range("A82").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A83").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A84").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A85").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A86").Value = 'SC-13,MA-4 (6)'
$ws.Range("A87").Value = 'AC-17 (2),MA-4 (6)'
$ws.Range("A88").Value = 'SC-13,MA-4 (6)'
$ws.Range("A92").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A93").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A94").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A95").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A96").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A97").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A98").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A99").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A100").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A101").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A102").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A103").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A104").Value = 'AC-11 (1),AC-11 b'
$ws.Range("A113").Value = 'AC-17 (2),SC-8,SC-13,MA-4 c'
$ws.Range("A114").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A115").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A116").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A117").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A122").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A123").Value = 'SC-28,SC-28 (1)'
$ws.Range("A133").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A134").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A135").Value = 'AC-6 (10),CM-6 b'
$ws.Range("A136").Value = 'AC-6 (10),AC-3 (4)'
$ws.Range("A137").Value = 'AC-6 (10),AC-3 (4)'
$ws.Range("A139").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A140").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A141").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A142").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A146").Value = 'IA-7,IA-5 (1) (c)'
$ws.Range("A152").Value = 'IA-7,CM-6 b'
$ws.Range("A153").Value = 'IA-7,CM-6 b'
$ws.Range("A154").Value = 'IA-7,CM-6 b'
$ws.Range("A157").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A158").Value = 'AU-3,AU-12 a,MA-4 (1) (a),AU-12 c'
$ws.Range("A159").Value = 'AU-3,AU-12 a,MA-4 (1) (a),AU-12 c'
$ws.Range("A160").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A161").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A162").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A163").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A164").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A171").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A175").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A176").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A177").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A178").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A179").Value = 'AU-3,MA-4 (1) (a),AU-3 (1)'
$ws.Range("A180").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A181").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A182").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A183").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A184").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A185").Value = 'MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A186").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A187").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A188").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A189").Value = 'AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A190").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A191").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A192").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A193").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A194").Value = 'AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A195").Value = 'AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)'
$ws.Range("A196").Value = 'IA-2 (1),IA-2 (2),IA-2 (4),IA-2 (3)'
$ws.Range("A197").Value = 'IA-2 (4),IA-2 (5),IA-2 (2),IA-2 (3),IA-2'
$ws.Range("A198").Value = 'IA-2 (4),IA-2 (5),IA-2 (2),IA-2 (3),IA-2'
$ws.Range("A206").Value = 'SC-8 (2),SC-8,SC-8 (1)'
$ws.Range("A207").Value = 'SC-8 (2),SC-8,SC-8 (1)'
$ws.Range("A208").Value = 'AC-18 (1),SC-8,SC-8 (1)'
$ws.Range("A215").Value = 'AU-12 a,CM-6 b'
$ws.Range("A216").Value = 'SC-5 (2),SC-5,CM-6 b'
$ws.Range("A219").Value = 'SI-16,CM-6 b'
$ws.Range("A220").Value = 'IA-8,AU-3 (1),IA-2'
$ws.Range("A232").Value = 'SI-16,SC-2,CM-6 b'
$ws.Range("A234").Value = 'SC-3,SI-16'
$ws.Range("A243").Value = 'CM-6 b,IA-5 (1) (b),IA-5 (1) (a)'
$ws.Range("A251").Value = 'SC-2,SC-4'
$ws.Range("A252").Value = 'SC-2,SC-4'
$ws.Range("A254").Value = 'IA-2 (11),IA-2 (12)'
$ws.Range("A258").Value = 'SC-3,SI-6 a'
$ws.Range("A275").Value = 'CM-7 a,CM-6 b'
$ws.Range("A294").Value = 'AU-5 (1),AU-5 a'
$ws.Range("A339").Value = 'IA-5 (1) (c),CM-6 b'
$ws.Range("A341").Value = 'IA-2 (2),CM-6 b'
$ws.Range("A342").Value = 'CM-5 (1),CM-6 b'
$ws.Range("A343").Value = 'CM-5 (1),CM-6 b'
$ws.Range("A351").Value = 'AC-17 (2),CM-6 b'
$ws.Range("A374").Value = 'AU-3,CM-6 b'
$ws.Range("A377").Value = 'SC-3,CM-6 b'
$ws.Range("A382").Value = 'AC-17 (9),CM-7 b,AC-17 (1),CM-6 b'
$ws.Range("A383").Value = 'CM-7 b,AC-17 (1),CM-6 b'
$ws.Range("A422").Value = 'SC-3,CM-6 b'
$ws.Range("A424").Value = 'SC-2,CM-6 b'
$ws.Range("A425").Value = 'SC-2,CM-6 b'
$ws.Range("A429").Value = 'SC-3,CM-6 b'
$ws.Range("A445").Value = 'SI-2 (2),CM-6 b'
$ws.Range("A454").Value = 'SI-2 (2),CM-6 b'
$ws.Range("A458").Value = 'MA-4 (7),SC-10,MA-4 e,AC-12'
$ws.Range("A461").Value = 'AC-11 a,SC-10'
$ws.Range("A462").Value = 'AC-17 (2),SC-8,SC-8 (1)'
$ws.Range("A479").Value = 'CM-7 b,IA-3'
$ws.Range("A497").Value = 'AU-4 (1),AU-4'
$ws.Range("A533").Value = 'AC-17 (2),SC-8'
